# Daily attendance processing - 2025-11-05 05:47:40
# Normalizes the "Recorded By" (column G) entries so that a human
# reviewer's address is listed before "System" (e.g. "System, x@y.com"
# becomes "x@y.com, System"). Entries already led by the automated
# "backup@backdoor.com" account, or with only a single recorder, are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $value = $cell.Value2

    if ($value -eq $null) {
        continue
    }

    $parts = $value -split ", "

    if ($parts.Length -eq 2 -and $parts[0] -ne "backup@backdoor.com") {
        $cell.Value = $parts[1] + ", " + $parts[0]
    }
}
